$wb = $excel.ActiveWorkbook

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
# Row 116
$ws.Cells.Item(116, 8).Value = 3001.2856
$ws.Cells.Item(116, 9).Value = 2252.5
$ws.Cells.Item(116, 11).Value = 2252.5
$ws.Cells.Item(116, 13).Value = 1189.5

# Row 132
$ws.Cells.Item(132, 8).Value = 5536.467
$ws.Cells.Item(132, 9).Value = 1448.56
$ws.Cells.Item(132, 10).Value = 25976
$ws.Cells.Item(132, 11).Value = 4345.68
$ws.Cells.Item(132, 12).Value = 77928
$ws.Cells.Item(132, 13).Value = -1815.68
$ws.Cells.Item(132, 14).Value = -82988

# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Cells.Item(32, 8).Value = 7098.628
$ws.Cells.Item(32, 9).Value = 3568.2466
$ws.Cells.Item(32, 10).Value = 26923.076
$ws.Cells.Item(32, 11).Value = 3568.2466
$ws.Cells.Item(32, 12).Value = 26923.076
$ws.Cells.Item(32, 13).Value = -3281.2466
$ws.Cells.Item(32, 14).Value = -27497.076

# Row 63
$ws.Cells.Item(63, 8).Value = 3992
$ws.Cells.Item(63, 9).Value = 2485
$ws.Cells.Item(63, 10).Value = 5714.2856
$ws.Cells.Item(63, 11).Value = 2485
$ws.Cells.Item(63, 12).Value = 5714.2856
$ws.Cells.Item(63, 13).Value = -1799
$ws.Cells.Item(63, 14).Value = -7086.2856

# Row 66
$ws.Cells.Item(66, 8).Value = 3992
$ws.Cells.Item(66, 9).Value = 2485
$ws.Cells.Item(66, 10).Value = 5714.2856
$ws.Cells.Item(66, 11).Value = 12425
$ws.Cells.Item(66, 12).Value = 28571.428
$ws.Cells.Item(66, 13).Value = -8993
$ws.Cells.Item(66, 14).Value = -35435.428

# Row 102
$ws.Cells.Item(102, 8).Value = 3559.8
$ws.Cells.Item(102, 9).Value = 2949.75
$ws.Cells.Item(102, 10).Value = 6000
$ws.Cells.Item(102, 11).Value = 2949.75
$ws.Cells.Item(102, 12).Value = 6000
$ws.Cells.Item(102, 13).Value = -1327.75
$ws.Cells.Item(102, 14).Value = -9244

# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
# Row 52
$ws.Cells.Item(52, 8).Value = 32470
$ws.Cells.Item(52, 10).Value = 32470
$ws.Cells.Item(52, 12).Value = 32470
$ws.Cells.Item(52, 14).Value = -32996

# Row 55
$ws.Cells.Item(55, 8).Value = 29296.572
$ws.Cells.Item(55, 10).Value = 29296.572
$ws.Cells.Item(55, 12).Value = 29296.572
$ws.Cells.Item(55, 14).Value = -29842.572

# Row 86
$ws.Cells.Item(86, 8).Value = 1851.4839
$ws.Cells.Item(86, 9).Value = 1632.7778
$ws.Cells.Item(86, 10).Value = 2154.3076
$ws.Cells.Item(86, 11).Value = 1632.7778
$ws.Cells.Item(86, 12).Value = 2154.3076
$ws.Cells.Item(86, 13).Value = -509.7778000000001
$ws.Cells.Item(86, 14).Value = -4400.3076

# Row 89
$ws.Cells.Item(89, 8).Value = 1851.4839
$ws.Cells.Item(89, 9).Value = 1632.7778
$ws.Cells.Item(89, 10).Value = 2154.3076
$ws.Cells.Item(89, 11).Value = 8163.889
$ws.Cells.Item(89, 12).Value = 10771.538
$ws.Cells.Item(89, 13).Value = -2547.889
$ws.Cells.Item(89, 14).Value = -22003.538

# Row 99
$ws.Cells.Item(99, 8).Value = 2119.875
$ws.Cells.Item(99, 9).Value = 1910
$ws.Cells.Item(99, 10).Value = 2749.5
$ws.Cells.Item(99, 11).Value = 1910
$ws.Cells.Item(99, 12).Value = 2749.5
$ws.Cells.Item(99, 13).Value = -412
$ws.Cells.Item(99, 14).Value = -5745.5

# Row 105
$ws.Cells.Item(105, 8).Value = 2100.0667
$ws.Cells.Item(105, 9).Value = 2053.8462
$ws.Cells.Item(105, 10).Value = 2400.5
$ws.Cells.Item(105, 11).Value = 2053.8462
$ws.Cells.Item(105, 12).Value = 2400.5
$ws.Cells.Item(105, 13).Value = -306.8462
$ws.Cells.Item(105, 14).Value = -5894.5

# Row 121
$ws.Cells.Item(121, 8).Value = 32470
$ws.Cells.Item(121, 10).Value = 32470
$ws.Cells.Item(121, 12).Value = 32470
$ws.Cells.Item(121, 14).Value = -35964

# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31, 8).Value = 2652.1035
$ws.Cells.Item(31, 9).Value = 1558.3
$ws.Cells.Item(31, 11).Value = 1558.3
$ws.Cells.Item(31, 13).Value = -1263.3

# Row 34
$ws.Cells.Item(34, 8).Value = 2652.1035
$ws.Cells.Item(34, 9).Value = 1558.3
$ws.Cells.Item(34, 11).Value = 1558.3
$ws.Cells.Item(34, 13).Value = -1356.3

# Row 105
$ws.Cells.Item(105, 8).Value = 1087.1428
$ws.Cells.Item(105, 9).Value = 1022
$ws.Cells.Item(105, 10).Value = 1250
$ws.Cells.Item(105, 11).Value = 1022
$ws.Cells.Item(105, 12).Value = 1250
$ws.Cells.Item(105, 13).Value = 725
$ws.Cells.Item(105, 14).Value = -4744

# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
# Row 100
$ws.Cells.Item(100, 8).Value = 3600
$ws.Cells.Item(100, 10).Value = 3600
$ws.Cells.Item(100, 12).Value = 10800
$ws.Cells.Item(100, 14).Value = -12422

# Row 107
$ws.Cells.Item(107, 8).Value = 814.9286
$ws.Cells.Item(107, 9).Value = 876.5833
$ws.Cells.Item(107, 10).Value = 445
$ws.Cells.Item(107, 11).Value = 2629.7499
$ws.Cells.Item(107, 12).Value = 1335
$ws.Cells.Item(107, 13).Value = -709.7498999999998
$ws.Cells.Item(107, 14).Value = -5175

# Row 115
$ws.Cells.Item(115, 8).Value = 2636.6667
$ws.Cells.Item(115, 9).Value = 1200
$ws.Cells.Item(115, 10).Value = 2924
$ws.Cells.Item(115, 11).Value = 3600
$ws.Cells.Item(115, 12).Value = 8772
$ws.Cells.Item(115, 13).Value = -2425
$ws.Cells.Item(115, 14).Value = -11122

# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
# Row 51
$ws.Cells.Item(51, 8).Value = 29200
$ws.Cells.Item(51, 10).Value = 29200
$ws.Cells.Item(51, 12).Value = 29200
$ws.Cells.Item(51, 14).Value = -30218

# Row 57
$ws.Cells.Item(57, 8).Value = 19119.8
$ws.Cells.Item(57, 10).Value = 19119.8
$ws.Cells.Item(57, 12).Value = 19119.8
$ws.Cells.Item(57, 14).Value = -20759.8

# Row 80
$ws.Cells.Item(80, 8).Value = 2290.9524
$ws.Cells.Item(80, 9).Value = 2363.3333
$ws.Cells.Item(80, 10).Value = 2194.4443
$ws.Cells.Item(80, 11).Value = 2363.3333
$ws.Cells.Item(80, 12).Value = 2194.4443
$ws.Cells.Item(80, 13).Value = -1365.3333
$ws.Cells.Item(80, 14).Value = -4190.4443

# Row 83
$ws.Cells.Item(83, 8).Value = 2290.9524
$ws.Cells.Item(83, 9).Value = 2363.3333
$ws.Cells.Item(83, 10).Value = 2194.4443
$ws.Cells.Item(83, 11).Value = 11816.6665
$ws.Cells.Item(83, 12).Value = 10972.2215
$ws.Cells.Item(83, 13).Value = -6824.666499999999
$ws.Cells.Item(83, 14).Value = -20956.2215

# Row 126
$ws.Cells.Item(126, 8).Value = 55556824
$ws.Cells.Item(126, 9).Value = 100001120
$ws.Cells.Item(126, 10).Value = 1460.625
$ws.Cells.Item(126, 11).Value = 300003360
$ws.Cells.Item(126, 12).Value = 4381.875
$ws.Cells.Item(126, 13).Value = -300000890
$ws.Cells.Item(126, 14).Value = -9321.875

# --- Sheet 7 (LTW) ---
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Cells.Item(7, 8).Value = 1352.174
$ws.Cells.Item(7, 9).Value = 1368.2632
$ws.Cells.Item(7, 10).Value = 1275.75
$ws.Cells.Item(7, 11).Value = 1368.2632
$ws.Cells.Item(7, 12).Value = 1275.75
$ws.Cells.Item(7, 13).Value = -1256.2632
$ws.Cells.Item(7, 14).Value = -1499.75

# Row 53
$ws.Cells.Item(53, 8).Value = 13224.75
$ws.Cells.Item(53, 9).Value = 9633
$ws.Cells.Item(53, 10).Value = 24000
$ws.Cells.Item(53, 11).Value = 9633
$ws.Cells.Item(53, 12).Value = 24000
$ws.Cells.Item(53, 13).Value = -9115
$ws.Cells.Item(53, 14).Value = -25036

# Row 68
$ws.Cells.Item(68, 8).Value = 8868
$ws.Cells.Item(68, 9).Value = 26350.5
$ws.Cells.Item(68, 10).Value = 1875
$ws.Cells.Item(68, 11).Value = 26350.5
$ws.Cells.Item(68, 12).Value = 1875
$ws.Cells.Item(68, 13).Value = -25601.5
$ws.Cells.Item(68, 14).Value = -3373

# Row 71
$ws.Cells.Item(71, 8).Value = 8868
$ws.Cells.Item(71, 9).Value = 26350.5
$ws.Cells.Item(71, 10).Value = 1875
$ws.Cells.Item(71, 11).Value = 131752.5
$ws.Cells.Item(71, 12).Value = 9375
$ws.Cells.Item(71, 13).Value = -128008.5
$ws.Cells.Item(71, 14).Value = -16863

# Row 93
$ws.Cells.Item(93, 8).Value = 70561.2
$ws.Cells.Item(93, 9).Value = 750
$ws.Cells.Item(93, 10).Value = 117102
$ws.Cells.Item(93, 11).Value = 750
$ws.Cells.Item(93, 12).Value = 117102
$ws.Cells.Item(93, 13).Value = 498
$ws.Cells.Item(93, 14).Value = -119598

# Row 100
$ws.Cells.Item(100, 8).Value = 1653.6666
$ws.Cells.Item(100, 9).Value = 1205.5
$ws.Cells.Item(100, 10).Value = 1877.75
$ws.Cells.Item(100, 11).Value = 1205.5
$ws.Cells.Item(100, 12).Value = 1877.75
$ws.Cells.Item(100, 13).Value = -664.5
$ws.Cells.Item(100, 14).Value = -2959.75

# Row 126
$ws.Cells.Item(126, 8).Value = 1352.174
$ws.Cells.Item(126, 9).Value = 1368.2632
$ws.Cells.Item(126, 10).Value = 1275.75
$ws.Cells.Item(126, 11).Value = 4104.7896
$ws.Cells.Item(126, 12).Value = 3827.25
$ws.Cells.Item(126, 13).Value = -1634.7896
$ws.Cells.Item(126, 14).Value = -8767.25

# --- Sheet 8 (WVR) ---
$ws = $wb.Worksheets.Item(8)
# Row 50
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).Value = $null

# Row 96
$ws.Cells.Item(96, 8).Value = 1249.8334
$ws.Cells.Item(96, 9).Value = 1176.5
$ws.Cells.Item(96, 10).Value = 1286.5
$ws.Cells.Item(96, 11).Value = 1176.5
$ws.Cells.Item(96, 12).Value = 1286.5
$ws.Cells.Item(96, 13).Value = 196.5
$ws.Cells.Item(96, 14).Value = -4032.5
